# Weekly update: insert 2 new observation rows (new week) right after the
# existing row for "Fruta, Feria Lagunitas de Puerto Montt - Mango" at row 144.
# This pushes every subsequent row down by two, matching the growth of the
# sheet dimension from A1:T168 to A1:T170.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows before the old row 145.
$ws.Rows.Item(145).Resize(2).Insert()

# --- New row 145: "Primera" quality entry for the new week ---
$ws.Cells.Item(145, 1).Value = 4
$ws.Cells.Item(145, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(145, 3).Value = "Los Lagos"
$ws.Cells.Item(145, 4).Value = 44617
$ws.Cells.Item(145, 5).Value = 10
$ws.Cells.Item(145, 6).Value = "Fruta"
$ws.Cells.Item(145, 7).Value = 100108
$ws.Cells.Item(145, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(145, 9).Value = 100108002
$ws.Cells.Item(145, 10).Value = "Mango"
$ws.Cells.Item(145, 11).Value = "Sin especificar"
$ws.Cells.Item(145, 12).Value = "Primera"
$ws.Cells.Item(145, 13).Value = 200
$ws.Cells.Item(145, 14).Value = 7000
$ws.Cells.Item(145, 15).Value = 7500
$ws.Cells.Item(145, 16).Value = 7250
$ws.Cells.Item(145, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(145, 18).Value = "Perú"
$ws.Cells.Item(145, 19).Value = 1812
$ws.Cells.Item(145, 20).Value = 4

# --- New row 146: "Segunda" quality entry for the same new week ---
$ws.Cells.Item(146, 1).Value = 4
$ws.Cells.Item(146, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(146, 3).Value = "Los Lagos"
$ws.Cells.Item(146, 4).Value = 44617
$ws.Cells.Item(146, 5).Value = 10
$ws.Cells.Item(146, 6).Value = "Fruta"
$ws.Cells.Item(146, 7).Value = 100108
$ws.Cells.Item(146, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(146, 9).Value = 100108002
$ws.Cells.Item(146, 10).Value = "Mango"
$ws.Cells.Item(146, 11).Value = "Sin especificar"
$ws.Cells.Item(146, 12).Value = "Segunda"
$ws.Cells.Item(146, 13).Value = 100
$ws.Cells.Item(146, 14).Value = 5000
$ws.Cells.Item(146, 15).Value = 5000
$ws.Cells.Item(146, 16).Value = 5000
$ws.Cells.Item(146, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(146, 18).Value = "Perú"
$ws.Cells.Item(146, 19).Value = 1250
$ws.Cells.Item(146, 20).Value = 4
